# daily auto push: 2026-01-15 22:35 UTC
# A new daily-scrape row for 2026/01/16 03:00 is inserted at the top of the
# "future" block (row 659), pushing every existing row from 659..700 down by
# one (to 660..701). The sheet's used range grows from D700 to D701.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 659:700 down to 660:701, leaving a blank row 659 for the new entry.
$ws.Rows(659).Insert()

# Column A holds dates formatted as plain text ("yyyy/mm/dd"), not real Excel
# dates. Flip the new cell to text first so Excel doesn't reinterpret the
# string as a date serial, then drop the explicit number format again so the
# cell matches the style-less cells around it.
$ws.Range("A659").NumberFormat = "@"
$ws.Range("A659").Value = "2026/01/16"
$ws.Range("A659").ClearFormats()

$ws.Range("B659").Value = "金"
$ws.Range("C659").Value = 3
$ws.Range("D659").Value = 201
